$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-detected as numbers (losing trailing zeros / exact text form),
# to match the original inline-string (text) representation.
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'

$ws.Range('D2').Value = '26.216.49'
$ws.Range('D3').Value = '1.582.88'
$ws.Range('E3').Value = '  -1.17%  '
$ws.Range('E4').Value = '  -0.41%  '
$ws.Range('D5').Value = '209.70'
$ws.Range('E5').Value = '  -0.93%  '
$ws.Range('E6').Value = '  -2.59%  '
$ws.Range('E7').Value = '  -0.34%  '
$ws.Range('E9').Value = '  -0.53%  '
$ws.Range('E10').Value = '  -1.06%  '
$ws.Range('D11').Value = '0.0847'
$ws.Range('E11').Value = '  +0.35%  '
$ws.Range('D12').Value = '1.804.91'
$ws.Range('E12').Value = '  -1.19%  '
$ws.Range('D13').Value = '1.591.67'
$ws.Range('E13').Value = '  -0.46%  '
$ws.Range('E14').Value = '  -0.16%  '
$ws.Range('E15').Value = '  -1.02%  '
$ws.Range('E16').Value = '  -0.61%  '
$ws.Range('D17').Value = '26.209.66'
$ws.Range('E17').Value = '  -1.82%  '
$ws.Range('E18').Value = '  -1.02%  '
$ws.Range('E19').Value = '  +1.05%  '
$ws.Range('E20').Value = '  -0.36%  '
$ws.Range('D21').Value = '206.19'
$ws.Range('D22').Value = '4.25'
$ws.Range('E22').Value = '  -0.75%  '
$ws.Range('E23').Value = '  -2.94%  '
$ws.Range('E24').Value = '  -0.85%  '
$ws.Range('D25').Value = '144.81'
$ws.Range('E25').Value = '  +0.33%  '
$ws.Range('E26').Value = '  -0.35%  '
$ws.Range('D27').Value = '7.03'
$ws.Range('E27').Value = '  -0.72%  '
$ws.Range('E28').Value = '  -1.01%  '
$ws.Range('D29').Value = '15.22'
$ws.Range('D31').Value = '1.15'
$ws.Range('E31').Value = '  -0.74%  '
$ws.Range('E32').Value = '  -1.40%  '
$ws.Range('E33').Value = '  -1.12%  '
$ws.Range('D34').Value = '1.283.42'
$ws.Range('E34').Value = '  -0.80%  '
$ws.Range('E35').Value = '  +8.01%  '
$ws.Range('E36').Value = '  -0.34%  '
$ws.Range('D37').Value = '0.605'
$ws.Range('E37').Value = '  +0.80%  '
$ws.Range('E38').Value = '  -1.17%  '
$ws.Range('E39').Value = '  -1.63%  '
$ws.Range('D40').Value = '0.814'
$ws.Range('E40').Value = '  -1.81%  '
$ws.Range('D41').Value = '5.58'
$ws.Range('E41').Value = '  +3.35%  '
$ws.Range('D42').Value = '0.769'
$ws.Range('E42').Value = '  -1.50%  '
$ws.Range('E43').Value = '  -2.85%  '
$ws.Range('D44').Value = '62.25'
$ws.Range('E44').Value = '  -1.21%  '
$ws.Range('D45').Value = '1.718.62'
$ws.Range('E45').Value = '  -1.25%  '
$ws.Range('E46').Value = '  -2.28%  '
$ws.Range('E47').Value = '  -0.41%  '
$ws.Range('D49').Value = '0.0507'
$ws.Range('E49').Value = '  -1.71%  '
$ws.Range('E50').Value = '  -0.16%  '
$ws.Range('E51').Value = '  +0.00%  '
